$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.676.65"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "3.017.13"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "510.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.58%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +2.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("E10").Value = "  +3.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.368"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.73%  "
$ws.Range("D12").Value = "3.540.85"
$ws.Range("E12").Value = "  +3.07%  "
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("E15").Value = "  +4.84%  "
$ws.Range("D16").Value = "56.667.04"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Value = "3.019.23"
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "333.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.73%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  +5.38%  "
$ws.Range("E24").Value = "  +4.68%  "
$ws.Range("D25").Value = "3.147.70"
$ws.Range("E25").Value = "  +3.08%  "
$ws.Range("E26").Value = "  +4.72%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "0.0₃0918"
$ws.Range("E28").Value = "  +9.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.04%  "
$ws.Range("E31").Value = "  +4.01%  "
$ws.Range("E32").Value = "  +3.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "153.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.68%  "
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "27.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.60%  "
$ws.Range("E38").Value = "  +3.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0662"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("D40").Value = "3.055.01"
$ws.Range("E40").Value = "  +3.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.59%  "
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("E43").Value = "  +5.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.655"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.17%  "
$ws.Range("D45").Value = "2.199.80"
$ws.Range("E45").Value = "  +3.99%  "
$ws.Range("E46").Value = "  +8.92%  "
$ws.Range("E47").Value = "  +2.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.927"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.34%  "
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0854"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.15%  "

